$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 7499.5
$ws.Range("J10").Value = 7499.5
$ws.Range("L10").Value = 7499.5
$ws.Range("N10").Value = -8085.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2936.125
$ws.Range("I40").Value = 1849
$ws.Range("K40").Value = 1849
$ws.Range("M40").Value = -1674

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2260.35
$ws.Range("J86").Value = 2107.2856
$ws.Range("L86").Value = 2107.2856
$ws.Range("N86").Value = -4353.2856

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1250.44
$ws.Range("J88").Value = 1505.2
$ws.Range("L88").Value = 1505.2
$ws.Range("N88").Value = -2317.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 2260.35
$ws.Range("J89").Value = 2107.2856
$ws.Range("L89").Value = 10536.428
$ws.Range("N89").Value = -21768.428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 1250.44
$ws.Range("J91").Value = 1505.2
$ws.Range("L91").Value = 1505.2
$ws.Range("N91").Value = -4313.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 4727.316
$ws.Range("I135").Value = 1824.6666
$ws.Range("J135").Value = 15612.25
$ws.Range("K135").Value = 16421.9994
$ws.Range("L135").Value = 140510.25
$ws.Range("M135").Value = -13886.9994
$ws.Range("N135").Value = -145580.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4027.1875
$ws.Range("I138").Value = 451.85715
$ws.Range("J138").Value = 5499.3823
$ws.Range("K138").Value = 1355.57145
$ws.Range("L138").Value = 16498.1469
$ws.Range("M138").Value = 3784.42855
$ws.Range("N138").Value = -26778.1469

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 197499.5
$ws.Range("J140").Value = 229999.33
$ws.Range("L140").Value = 229999.33
$ws.Range("N140").Value = -240359.33

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3950.5
$ws.Range("I61").Value = 4168.5654
$ws.Range("J61").Value = 3393.2222
$ws.Range("K61").Value = 4168.5654
$ws.Range("L61").Value = 3393.2222
$ws.Range("M61").Value = -3956.5654
$ws.Range("N61").Value = -3817.2222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4025.818
$ws.Range("I74").Value = 5171.643
$ws.Range("J74").Value = 2020.625
$ws.Range("K74").Value = 5171.643
$ws.Range("L74").Value = 2020.625
$ws.Range("M74").Value = -4297.643
$ws.Range("N74").Value = -3768.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4025.818
$ws.Range("I77").Value = 5171.643
$ws.Range("J77").Value = 2020.625
$ws.Range("K77").Value = 25858.215
$ws.Range("L77").Value = 10103.125
$ws.Range("M77").Value = -21490.215
$ws.Range("N77").Value = -18839.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3950.5
$ws.Range("I136").Value = 4168.5654
$ws.Range("J136").Value = 3393.2222
$ws.Range("K136").Value = 12505.6962
$ws.Range("L136").Value = 10179.6666
$ws.Range("M136").Value = -9955.696200000002
$ws.Range("N136").Value = -15279.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 353.95456
$ws.Range("I5").Value = 425.8
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 425.8
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = -312.8
$ws.Range("N5").Value = -426

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 15402.479
$ws.Range("I86").Value = 12632.454
$ws.Range("J86").Value = 17941.666
$ws.Range("K86").Value = 12632.454
$ws.Range("L86").Value = 17941.666
$ws.Range("M86").Value = -11509.454
$ws.Range("N86").Value = -20187.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 15402.479
$ws.Range("I89").Value = 12632.454
$ws.Range("J89").Value = 17941.666
$ws.Range("K89").Value = 63162.27
$ws.Range("L89").Value = 89708.33
$ws.Range("M89").Value = -57546.27
$ws.Range("N89").Value = -100940.33

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 32427.229
$ws.Range("I94").Value = 851.1053000000001
$ws.Range("K94").Value = 851.1053000000001
$ws.Range("M94").Value = -400.1053000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1255000.1
$ws.Range("I134").Value = 1518818.4
$ws.Range("J134").Value = 11285.571
$ws.Range("K134").Value = 4556455.199999999
$ws.Range("L134").Value = 33856.713
$ws.Range("M134").Value = -4553920.199999999
$ws.Range("N134").Value = -38926.713

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 75000
$ws.Range("J64").Value = 75000
$ws.Range("L64").Value = 75000
$ws.Range("N64").Value = -75496

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H67").Value = 75000
$ws.Range("J67").Value = 75000
$ws.Range("L67").Value = 75000
$ws.Range("N67").Value = -76716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 12131.615
$ws.Range("I86").Value = 10159.2
$ws.Range("K86").Value = 10159.2
$ws.Range("M86").Value = -9036.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 12131.615
$ws.Range("I89").Value = 10159.2
$ws.Range("K89").Value = 50796
$ws.Range("M89").Value = -45180

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 125003210
$ws.Range("I105").Value = 250002720
$ws.Range("J105").Value = 3699.75
$ws.Range("K105").Value = 250002720
$ws.Range("L105").Value = 3699.75
$ws.Range("M105").Value = -250000973
$ws.Range("N105").Value = -7193.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 90000
$ws.Range("N101").Value = -94868

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 5504.12
$ws.Range("I107").Value = 287.25
$ws.Range("J107").Value = 6497.8096
$ws.Range("K107").Value = 861.75
$ws.Range("L107").Value = 19493.4288
$ws.Range("M107").Value = 1058.25
$ws.Range("N107").Value = -23333.4288

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3714.0908
$ws.Range("I132").Value = 1946.5
$ws.Range("J132").Value = 4106.8887
$ws.Range("K132").Value = 17518.5
$ws.Range("L132").Value = 36961.99830000001
$ws.Range("M132").Value = -14988.5
$ws.Range("N132").Value = -42021.99830000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 1939.6471
$ws.Range("I137").Value = 1152.8334
$ws.Range("J137").Value = 2368.818
$ws.Range("K137").Value = 3458.5002
$ws.Range("L137").Value = 7106.454000000001
$ws.Range("M137").Value = 1641.4998
$ws.Range("N137").Value = -17306.454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 20000
$ws.Range("I48").Value = 20000
$ws.Range("K48").Value = 20000
$ws.Range("M48").Value = -19515

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 23033.5
$ws.Range("J49").Value = 23033.5
$ws.Range("L49").Value = 23033.5
$ws.Range("N49").Value = -23401.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13899.2
$ws.Range("I70").Value = 7674
$ws.Range("K70").Value = 7674
$ws.Range("M70").Value = -7404

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 13899.2
$ws.Range("I73").Value = 7674
$ws.Range("K73").Value = 7674
$ws.Range("M73").Value = -6738

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4600.491
$ws.Range("I102").Value = 3189.081
$ws.Range("K102").Value = 3189.081
$ws.Range("M102").Value = -1567.081

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 631.93335
$ws.Range("I107").Value = 558.8889
$ws.Range("J107").Value = 741.5
$ws.Range("K107").Value = 558.8889
$ws.Range("L107").Value = 741.5
$ws.Range("M107").Value = 1361.1111
$ws.Range("N107").Value = -4581.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 16671449
$ws.Range("I126").Value = 26318908
$ws.Range("K126").Value = 78956724
$ws.Range("M126").Value = -78954254

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4931.7617
$ws.Range("I40").Value = 4248.5713
$ws.Range("K40").Value = 4248.5713
$ws.Range("M40").Value = -4112.5713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 12823645
$ws.Range("I136").Value = 16132128
$ws.Range("J136").Value = 3274.75
$ws.Range("K136").Value = 48396384
$ws.Range("L136").Value = 9824.25
$ws.Range("M136").Value = -48393834
$ws.Range("N136").Value = -14924.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1895.2778
$ws.Range("I122").Value = 1626.3846
$ws.Range("J122").Value = 2594.4
$ws.Range("K122").Value = 4879.1538
$ws.Range("L122").Value = 7783.200000000001
$ws.Range("M122").Value = -2429.1538
$ws.Range("N122").Value = -12683.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4522.5713
$ws.Range("J132").Value = 9386.286
$ws.Range("L132").Value = 28158.858
$ws.Range("N132").Value = -33218.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 9440986
$ws.Range("I136").Value = 12821738
$ws.Range("K136").Value = 38465214
$ws.Range("M136").Value = -38462664
